# Add a new weekly price record for "Acelga" (Vega Modelo de Temuco) as row 186,
# shifting the existing rows 186-214 down to 187-215.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 186 (pushes old rows 186..214 -> 187..215).
$ws.Rows.Item(186).Insert()

# Populate the new row 186 with the new record's data.
$ws.Range("A186").Value = 10
$ws.Range("B186").Value = "Vega Modelo de Temuco"
$ws.Range("C186").Value = "La Araucanía"
$ws.Range("D186").Value = 44504
$ws.Range("E186").Value = 9
$ws.Range("F186").Value = 100112009
$ws.Range("G186").Value = "Acelga"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 115
$ws.Range("K186").Value = 8000
$ws.Range("L186").Value = 9000
$ws.Range("M186").Value = 8478
$ws.Range("N186").Value = "$/docena de atados (12 kilos)"
$ws.Range("O186").Value = "Provincia de Cautín"
$ws.Range("P186").Value = 706
$ws.Range("Q186").Value = 12
$ws.Range("R186").Value = "Hortaliza"
